$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Grab the "intro" row formatting (ht=25, big Arial font) that currently
#    lives on rows 2 & 3 and stamp it onto two new blank rows (29 & 30)
#    at the bottom of the sheet before we touch rows 2/3 at all.
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A29").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(29).RowHeight = 25
$ws.Rows.Item(30).RowHeight = 25

# ---------------------------------------------------------------------------
# 2) Strip that special formatting back off rows 2 & 3 - in the new layout
#    every item row (2-24) uses plain/default formatting.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).ClearFormats() | Out-Null
$ws.Rows.Item(2).AutoFit() | Out-Null
$ws.Rows.Item(3).ClearFormats() | Out-Null
$ws.Rows.Item(3).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 3) Re-write the 23 survey items (rows 2-24) in their corrected, fully
#    numeric order (01-23), fixing the stray punctuation/spacing typos
#    along the way.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = '01. I wear blue-filtering, orange-tinted, and/or red-tinted glasses indoors during the day.'
$ws.Range("A3").Value = '02. I wear blue-filtering, orange-tinted, and/or red-tinted glasses outdoors during the day.'
$ws.Range("A4").Value = '03. I wear blue-filtering, orange-tinted, and/or red-tinted glasses within 1 hour before attempting to fall asleep.'
$ws.Range("A5").Value = '04. I spend 30 minutes or less per day (in total) outside.'
$ws.Range("A6").Value = '05. I spend between 30 minutes and 1 hour per day (in total) outside.'
$ws.Range("A7").Value = '06. I spend between 1 and 3 hours per day (in total) outside.'
$ws.Range("A8").Value = '07. I spend more than 3 hours per day (in total) outside.'
$ws.Range("A9").Value = '08. I spend as much time outside as possible.'
$ws.Range("A10").Value = '09. I go for a walk or exercise outside within 2 hours after waking up.'
$ws.Range("A11").Value = '10. I use my mobile phone within 1 hour before attempting to fall asleep.'
$ws.Range("A12").Value = '11. I look at my mobile phone screen immediately after waking up.'
$ws.Range("A13").Value = '12. I check my phone when I wake up at night.'
$ws.Range("A14").Value = '13. I look at my smartwatch within 1 hour before attempting to fall asleep'
$ws.Range("A15").Value = '14. I look at my smartwatch when I wake up at night.'
$ws.Range("A16").Value = '15. I dim my mobile phone screen within 1 hour before attempting to fall asleep.'
$ws.Range("A17").Value = '16. I use a blue-filter app on my computer screen within 1 hour before attempting to fall asleep.'
$ws.Range("A18").Value = '17. I use as little light as possible when I get up during the night.'
$ws.Range("A19").Value = '18. I dim my computer screen within 1 hour before attempting to fall asleep.'
$ws.Range("A20").Value = '19. I use tunable lights to create a healthy light environment.'
$ws.Range("A21").Value = '20. I use LEDs to create a healthy light environment.'
$ws.Range("A22").Value = '21. I use a desk lamp when I do focused work.'
$ws.Range("A23").Value = '22. I use an alarm with a dawn simulation light.'
$ws.Range("A24").Value = '23. I turn on the lights immediately after waking up.'

# ---------------------------------------------------------------------------
# 4) Header row - bold the whole row; centre the first header cell.
# ---------------------------------------------------------------------------
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("B1:F1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 5) Column width / view tweaks.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 136.83203125

$excel.ActiveWindow.Zoom = 170
$ws.Range("A29").Select() | Out-Null

Write-Output "done"
